$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new "Save" column header in H1, copying the formatting (bold, border,
# alignment) used by the other header cells such as G1 ("sum").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for the existing data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
